$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '38.741.40'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '2.096.79'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '226.67'
$ws.Range('E5').Value = '  -0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '61.71'
$ws.Range('E7').Value = '  +2.45%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '0.0837'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '15.69'
$ws.Range('E12').Value = '  +5.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '2.410.31'
$ws.Range('E13').Value = '  +0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '21.85'
$ws.Range('E14').Value = '  -2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '0.799'
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '5.48'
$ws.Range('E16').Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '2.097.86'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '38.740.42'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '71.41'
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '6.05'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0841'
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '227.04'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '2.35'
$ws.Range('E24').Value = '  -3.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '2.30'
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '9.64'
$ws.Range('E26').Value = '  +2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '170.78'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('E29').Value = '  +2.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '19.27'
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '2.55'
$ws.Range('E31').Value = '  +10.01%  '
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('E33').Value = '  +13.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '4.71'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '0.0613'
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range('D37').Value = '2.36'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '3.49'
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '17.88'
$ws.Range('E40').Value = '  -1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '0.0226'
$ws.Range('E41').Value = '  +2.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range('D42').Value = '101.10'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '1.520.41'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('E44').Value = '  +6.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '2.81'
$ws.Range('E45').Value = '  +0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '7.77'
$ws.Range('E46').Value = '  +1.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '0.0909'
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('E48').Value = '  +4.79%  '
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '2.95'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '2.297.19'
$ws.Range('E51').Value = '  +0.91%  '
